# v2.6 - Added decoupled suspension, four-wheel steering, scripts to generate GGV diagram
#
# This script:
#  1. Updates the aero coefficient (CD @ H8) on the Sedan_Hamba sheet.
#  2. Adds a new "FSAE_Achilles" worksheet (cloned from the existing
#     per-vehicle Aero template) at the end of the workbook and fills in
#     its aero coefficients.
#  3. Makes Sedan_Hamba the active sheet (matching the new workbookView).

$wb = $excel.ActiveWorkbook

# --- 1. Sedan_Hamba: CD coefficient revised -----------------------------
$sedanHamba = $wb.Worksheets.Item("Sedan_Hamba")
$sedanHamba.Range("H8").Value = 1.98
$sedanHamba.Range("H8").NumberFormat = "0.00"

# --- 2. New vehicle sheet: FSAE_Achilles ---------------------------------
$template = $wb.Worksheets.Item("Trailer_Kumanzi")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "FSAE_Achilles"

# Instance name (row 3, column H)
$newSheet.Range("H3").Value = "FSAE_Achilles"

# Aero coefficients for the new vehicle
$newSheet.Range("H5").Value = -2.5
$newSheet.Range("H6").Value = 1
$newSheet.Range("H8").Value = 1.2
$newSheet.Range("F9").Value = -0.8
$newSheet.Range("G9").Value = 0
$newSheet.Range("H9").Value = 0.6

[void]$newSheet.Range("H3:H9").Select()
[void]$newSheet.Range("G12").Select()

# --- 3. Sedan_Hamba becomes the active sheet -----------------------------
[void]$sedanHamba.Activate()
[void]$sedanHamba.Select()
